$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: set Brand (column B) for rows 2-21 in row order
# Step 2: set Item Name (column D) for rows 2-21 in row order
# Step 3: set UOM (column E) for rows 2-21 in row order
# Step 4: set BSL NO (A) and ISL NO (C) numeric values
# Step 5: delete row 22

# --- Column B (BRAND) ---
$ws.Range('B2').Value = 'Aldorin'
$ws.Range('B3').Value = 'Cardoneb'
$ws.Range('B4').Value = 'Cardoneb'
$ws.Range('B5').Value = 'Cardovan'
$ws.Range('B6').Value = 'Cardovan'
$ws.Range('B7').Value = 'Cardovan'
$ws.Range('B8').Value = 'Cardovan'
$ws.Range('B9').Value = 'Dialon'
$ws.Range('B10').Value = 'GLIKAZID'
$ws.Range('B11').Value = 'Irbes'
$ws.Range('B12').Value = 'Ligazid'
$ws.Range('B13').Value = 'Ligazid'
$ws.Range('B14').Value = 'Ligazid'
$ws.Range('B15').Value = 'Lipicon'
$ws.Range('B16').Value = 'Lipicon'
$ws.Range('B17').Value = 'Lipicon'
$ws.Range('B18').Value = 'Lipicon'
$ws.Range('B19').Value = 'Pivasta'
$ws.Range('B20').Value = 'Sitazid'
$ws.Range('B21').Value = 'Sitazid'

# --- Column D (Item Name) ---
$ws.Range('D2').Value = 'Aldorin 50mg Tablet - 24''s'
$ws.Range('D3').Value = 'Cardoneb 5 FC Tablet 30''s'
$ws.Range('D4').Value = 'Cardoneb 2.5 FC Tablet 30''s'
$ws.Range('D5').Value = 'Cardovan Plus 80/12.5 Tablet 30''s'
$ws.Range('D6').Value = 'Cardovan 160mg Tablet 30''s'
$ws.Range('D7').Value = 'Cardovan 80mg Tablet 30''s'
$ws.Range('D8').Value = 'Cardovan Plus 160/12.5 Tablet 30''s'
$ws.Range('D9').Value = 'Dialon 4mg Tablet'
$ws.Range('D10').Value = 'Glikazid 80mg Tablet 30''s'
$ws.Range('D11').Value = 'Irbes 75mg Tablet'
$ws.Range('D12').Value = 'Ligazid 5mg Tablet 20''s'
$ws.Range('D13').Value = 'Ligazid 5mg Tablet 10''s'
$ws.Range('D14').Value = 'Ligazid M 2.5/500'
$ws.Range('D15').Value = 'Lipicon 20mg Tablet - 20''s'
$ws.Range('D16').Value = 'Lipicon 10mg Tablet Container 30''s'
$ws.Range('D17').Value = 'Lipicon 40mg Tablet - 10''s'
$ws.Range('D18').Value = 'Lipicon 10mg Tablet - 40''s'
$ws.Range('D19').Value = 'Pivasta 2mg Tablet 20''s'
$ws.Range('D20').Value = 'Sitazid 50mg Tablet 20''s'
$ws.Range('D21').Value = 'Sitazid 100mg Tablet 10''s'

# --- Column E (UOM) ---
$ws.Range('E2').Value = '24''s'
$ws.Range('E3').Value = '30''s'
$ws.Range('E4').Value = '30''s'
$ws.Range('E5').Value = '30''s'
$ws.Range('E6').Value = '30''s'
$ws.Range('E7').Value = '30''s'
$ws.Range('E8').Value = '30''s'
$ws.Range('E9').Value = '20''S'
$ws.Range('E10').Value = '30''s'
$ws.Range('E11').Value = '50 ''s'
$ws.Range('E12').Value = '20''s'
$ws.Range('E13').Value = '10''s'
$ws.Range('E14').Value = '20''s'
$ws.Range('E15').Value = '20 ''s'
$ws.Range('E16').Value = '30''s'
$ws.Range('E17').Value = '10 ''s'
$ws.Range('E18').Value = '40 ''s'
$ws.Range('E19').Value = '20''s'
$ws.Range('E20').Value = '20''s'
$ws.Range('E21').Value = '10''s'

# --- Column A (BSL NO) and C (ISL NO) ---
$ws.Range('A2').Value = 4
$ws.Range('C2').Value = 1
$ws.Range('A3').Value = 17
$ws.Range('C3').Value = 2
$ws.Range('A4').Value = 17
$ws.Range('C4').Value = 3
$ws.Range('A5').Value = 18
$ws.Range('C5').Value = 4
$ws.Range('A6').Value = 18
$ws.Range('C6').Value = 5
$ws.Range('A7').Value = 18
$ws.Range('C7').Value = 6
$ws.Range('A8').Value = 18
$ws.Range('C8').Value = 7
$ws.Range('A9').Value = 26
$ws.Range('C9').Value = 8
$ws.Range('A10').Value = 52
$ws.Range('C10').Value = 9
$ws.Range('A11').Value = 57
$ws.Range('C11').Value = 10
$ws.Range('A12').Value = 68
$ws.Range('C12').Value = 11
$ws.Range('A13').Value = 68
$ws.Range('C13').Value = 12
$ws.Range('A14').Value = 68
$ws.Range('C14').Value = 13
$ws.Range('A15').Value = 70
$ws.Range('C15').Value = 14
$ws.Range('A16').Value = 70
$ws.Range('C16').Value = 15
$ws.Range('A17').Value = 70
$ws.Range('C17').Value = 16
$ws.Range('A18').Value = 70
$ws.Range('C18').Value = 17
$ws.Range('A19').Value = 104
$ws.Range('C19').Value = 18
$ws.Range('A20').Value = 123
$ws.Range('C20').Value = 19
$ws.Range('A21').Value = 123
$ws.Range('C21').Value = 20

# --- Remove extra row 22 ---
$ws.Rows.Item(22).Delete()
